$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.224.41'
Set-TextValue $ws.Range('E2') '  +0.15%  '
Set-TextValue $ws.Range('D3') '1.832.98'
Set-TextValue $ws.Range('E3') '  -0.46%  '
Set-TextValue $ws.Range('D4') '0.9989'
Set-TextValue $ws.Range('E4') '  -0.37%  '
Set-TextValue $ws.Range('D5') '242.43'
Set-TextValue $ws.Range('E5') '  -0.79%  '
Set-TextValue $ws.Range('D6') '0.6213'
Set-TextValue $ws.Range('E6') '  -0.24%  '
Set-TextValue $ws.Range('D7') '1.000'
Set-TextValue $ws.Range('E7') '  -0.39%  '
Set-TextValue $ws.Range('D8') '0.07369'
Set-TextValue $ws.Range('E8') '  -1.84%  '
Set-TextValue $ws.Range('D9') '0.2911'
Set-TextValue $ws.Range('E9') '  -0.81%  '
Set-TextValue $ws.Range('D10') '23.24'
Set-TextValue $ws.Range('E10') '  -0.14%  '
Set-TextValue $ws.Range('D11') '0.07670'
Set-TextValue $ws.Range('E11') '  -0.69%  '
Set-TextValue $ws.Range('D12') '1.823.05'
Set-TextValue $ws.Range('E12') '  -1.82%  '
Set-TextValue $ws.Range('D13') '4.974'
Set-TextValue $ws.Range('E13') '  -0.86%  '
Set-TextValue $ws.Range('D14') '0.6691'
Set-TextValue $ws.Range('E14') '  -0.82%  '
Set-TextValue $ws.Range('D15') '82.68'
Set-TextValue $ws.Range('E15') '  -0.46%  '
Set-TextValue $ws.Range('D16') '0.000008964'
Set-TextValue $ws.Range('E16') '  -3.57%  '
Set-TextValue $ws.Range('D17') '5.869'
Set-TextValue $ws.Range('E17') '  -1.58%  '
Set-TextValue $ws.Range('D18') '29.202.41'
Set-TextValue $ws.Range('E18') '  -0.03%  '
Set-TextValue $ws.Range('D19') '2.072.85'
Set-TextValue $ws.Range('E19') '  -2.87%  '
Set-TextValue $ws.Range('D20') '236.08'
Set-TextValue $ws.Range('E20') '  +1.82%  '
Set-TextValue $ws.Range('E21') '  -1.54%  '
Set-TextValue $ws.Range('D22') '0.9999'
Set-TextValue $ws.Range('E22') '  -0.48%  '
Set-TextValue $ws.Range('D23') '7.358'
Set-TextValue $ws.Range('E23') '  +2.25%  '
Set-TextValue $ws.Range('D24') '1.001'
Set-TextValue $ws.Range('E24') '  -0.27%  '
Set-TextValue $ws.Range('E25') '  -1.48%  '
Set-TextValue $ws.Range('D26') '0.1403'
Set-TextValue $ws.Range('E26') '  +0.81%  '
Set-TextValue $ws.Range('D27') '8.541'
Set-TextValue $ws.Range('E27') '  -0.07%  '
Set-TextValue $ws.Range('D28') '17.65'
Set-TextValue $ws.Range('E28') '  -1.49%  '
Set-TextValue $ws.Range('D29') '1.488'
Set-TextValue $ws.Range('E29') '  -1.22%  '
Set-TextValue $ws.Range('D30') '0.05752'
Set-TextValue $ws.Range('E30') '  +2.98%  '
Set-TextValue $ws.Range('D31') '4.106'
Set-TextValue $ws.Range('E31') '  -1.09%  '
Set-TextValue $ws.Range('D32') '4.088'
Set-TextValue $ws.Range('E32') '  -2.25%  '
Set-TextValue $ws.Range('D33') '1.212'
Set-TextValue $ws.Range('E33') '  +0.54%  '
Set-TextValue $ws.Range('D34') '1.866'
Set-TextValue $ws.Range('E34') '  +0.98%  '
Set-TextValue $ws.Range('D35') '0.7284'
Set-TextValue $ws.Range('E35') '  -3.18%  '
Set-TextValue $ws.Range('D36') '1.142'
Set-TextValue $ws.Range('E36') '  -0.16%  '
Set-TextValue $ws.Range('D37') '2.602'
Set-TextValue $ws.Range('E37') '  -2.44%  '
Set-TextValue $ws.Range('D38') '2.860'
Set-TextValue $ws.Range('E38') '  +3.00%  '
Set-TextValue $ws.Range('D39') '1.223.84'
Set-TextValue $ws.Range('D40') '0.01755'
Set-TextValue $ws.Range('E40') '  -1.72%  '
Set-TextValue $ws.Range('D41') '6.276'
Set-TextValue $ws.Range('E41') '  -3.70%  '
Set-TextValue $ws.Range('D42') '0.9063'
Set-TextValue $ws.Range('E42') '  +0.66%  '
Set-TextValue $ws.Range('D43') '1.001'
Set-TextValue $ws.Range('E43') '  -0.29%  '
Set-TextValue $ws.Range('D44') '101.54'
Set-TextValue $ws.Range('E44') '  -0.58%  '
Set-TextValue $ws.Range('D45') '1.972.88'
Set-TextValue $ws.Range('E45') '  -2.75%  '
Set-TextValue $ws.Range('D46') '65.36'
Set-TextValue $ws.Range('E46') '  -1.01%  '
Set-TextValue $ws.Range('D47') '0.5042'
Set-TextValue $ws.Range('E47') '  -1.27%  '
Set-TextValue $ws.Range('E48') '  -2.97%  '
Set-TextValue $ws.Range('D49') '9.147'
Set-TextValue $ws.Range('E49') '  -0.12%  '
Set-TextValue $ws.Range('D50') '0.4025'
Set-TextValue $ws.Range('E50') '  -1.56%  '
Set-TextValue $ws.Range('D51') '0.1135'
Set-TextValue $ws.Range('E51') '  +3.02%  '
